$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.753.83"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = "'1.905.76"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.57%  '
$ws.Range('D4').Value = "'0.9978"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('D5').Value = "'312.47"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('D6').Value = "'0.9985"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').Value = "'0.5224"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +7.12%  '
$ws.Range('D8').Value = "'0.3785"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('D9').Value = "'0.07237"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.26%  '
$ws.Range('D10').Value = "'21.34"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +3.90%  '
$ws.Range('D11').Value = "'0.9077"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.76%  '
$ws.Range('D12').Value = "'0.07630"
$ws.Range('D12').ClearFormats()
$ws.Range('D13').Value = "'1.906.70"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('D14').Value = "'5.453"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').Value = "'92.23"
$ws.Range('D15').ClearFormats()
$ws.Range('D16').Value = "'0.9979"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.43%  '
$ws.Range('D17').Value = "'0.000008700"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.01%  '
$ws.Range('D18').Value = "'0.9986"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('D19').Value = "'27.785.80"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('D21').Value = "'5.150"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.48%  '
$ws.Range('D22').Value = "'2.148.36"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.51%  '
$ws.Range('D23').Value = "'10.85"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.01%  '
$ws.Range('D24').Value = "'6.603"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('D25').Value = "'153.33"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').Value = "'1.869"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.92%  '
$ws.Range('E27').Value = '  +0.82%  '
$ws.Range('D28').Value = "'18.32"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.27%  '
$ws.Range('D29').Value = "'114.50"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.08%  '
$ws.Range('D30').Value = "'4.848"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.06%  '
$ws.Range('D31').Value = "'0.08995"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.85%  '
$ws.Range('D32').Value = "'4.879"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +5.10%  '
$ws.Range('D33').Value = "'3.180"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('D34').Value = "'1.232"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.92%  '
$ws.Range('D35').Value = "'0.7766"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.41%  '
$ws.Range('D36').Value = "'2.628"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +4.04%  '
$ws.Range('D37').Value = "'0.02087"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.52%  '
$ws.Range('D38').Value = "'3.067"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.82%  '
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('D40').Value = "'0.5534"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.91%  '
$ws.Range('D41').Value = "'0.05283"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('D42').Value = "'6.690"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.24%  '
$ws.Range('D43').Value = "'114.82"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.51%  '
$ws.Range('D44').Value = "'8.543"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.54%  '
$ws.Range('D45').Value = "'0.1512"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('D46').Value = "'0.4818"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.42%  '
$ws.Range('D47').Value = "'10.47"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.08%  '
$ws.Range('D48').Value = "'0.9985"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.22%  '
$ws.Range('D49').Value = "'1.622"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('D50').Value = "'66.89"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('D51').Value = "'0.05993"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.92%  '
